$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-08-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-17 Saturday", 2)

# Update the multiplication problems in the table, cell by cell, so that
# duplicate texts (e.g. "144×8=1152" which appears twice with different
# replacements) are each replaced independently while preserving the
# existing run formatting of every cell.
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Text = "178×6=1068" },
    @{ Row = 1;  Col = 2; Text = "754×9=6786" },
    @{ Row = 1;  Col = 3; Text = "240×8=1920" },
    @{ Row = 1;  Col = 4; Text = "617×7=4319" },
    @{ Row = 1;  Col = 5; Text = "807×9=7263" },

    @{ Row = 5;  Col = 1; Text = "926×6=5556" },
    @{ Row = 5;  Col = 2; Text = "811×7=5677" },
    @{ Row = 5;  Col = 3; Text = "169×2=338" },
    @{ Row = 5;  Col = 4; Text = "186×8=1488" },
    @{ Row = 5;  Col = 5; Text = "449×7=3143" },

    @{ Row = 10; Col = 1; Text = "614×2=1228" },
    @{ Row = 10; Col = 2; Text = "300×9=2700" },
    @{ Row = 10; Col = 3; Text = "964×5=4820" },
    @{ Row = 10; Col = 4; Text = "969×5=4845" },
    @{ Row = 10; Col = 5; Text = "456×7=3192" },

    @{ Row = 15; Col = 1; Text = "242×8=1936" },
    @{ Row = 15; Col = 2; Text = "969×3=2907" },
    @{ Row = 15; Col = 3; Text = "399×4=1596" },
    @{ Row = 15; Col = 4; Text = "731×5=3655" },
    @{ Row = 15; Col = 5; Text = "139×3=417" },

    @{ Row = 20; Col = 1; Text = "806×2=1612" },
    @{ Row = 20; Col = 2; Text = "200×9=1800" },
    @{ Row = 20; Col = 3; Text = "778×3=2334" },
    @{ Row = 20; Col = 4; Text = "944×5=4720" },
    @{ Row = 20; Col = 5; Text = "432×2=864" }
)

foreach ($r in $replacements) {
    $cell = $tbl.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.Text
}
